# Generate Report for Handback
# Updates the "generate date" / handoff / handback timestamp cells on the
# Overview, zh-cn and de-de sheets to reflect freshly generated values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 8952ca72-...md
# (this text is shared with the de-de "Correspond Handoff Datetime" cell below)
$wsOverview.Range("G2").Value = "2016-10-20 00:48:04"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 8952ca72-...md
$wsZhCn.Range("H2").Value = "2016-10-20 00:47:52"
$wsZhCn.Range("K2").Value = "2016-10-20 00:48:34"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 8952ca72-...md
$wsDeDe.Range("H2").Value = "2016-10-20 00:48:04"
$wsDeDe.Range("K2").Value = "2016-10-20 00:48:52"
